# Updates the crypto price/volume table (Price = column D, Volume(1h) = column E)
# with the latest scraped values. Cells in column D that look like plain
# decimal numbers (e.g. "589.27") are forced to Text format first so Excel
# doesn't auto-convert them to numbers and strip meaningful formatting
# (trailing zeros, etc.) - matching the original inline-string text values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "66.885.66"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "2.597.48"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.27"
$ws.Range("E5").Value = "  -2.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.48"
$ws.Range("E6").Value = "  -3.18%  "
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("D9").Value = "2.597.02"
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.124"
$ws.Range("E10").Value = "  -2.49%  "
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("E13").Value = "  -2.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.18"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").Value = "3.067.82"
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("E16").Value = "  -4.87%  "
$ws.Range("D17").Value = "66.809.20"
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("D18").Value = "2.594.58"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "362.06"
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.00"
$ws.Range("E20").Value = "  -1.39%  "
$ws.Range("E21").Value = "  -4.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.29"
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.83"
$ws.Range("E23").Value = "  -1.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.04"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.35"
$ws.Range("E25").Value = "  +2.91%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("D28").Value = "2.727.14"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "580.95"
$ws.Range("E30").Value = "  +1.02%  "
$ws.Range("D31").Value = "0.0₃0980"
$ws.Range("E31").Value = "  -5.68%  "
$ws.Range("E32").Value = "  -4.73%  "
$ws.Range("E33").Value = "  -3.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.80"
$ws.Range("E34").Value = "  -3.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  -5.62%  "
$ws.Range("E37").Value = "  -2.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.11"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.93"
$ws.Range("E39").Value = "  -2.14%  "
$ws.Range("E40").Value = "  -1.20%  "
$ws.Range("E41").Value = "  -0.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.19"
$ws.Range("E42").Value = "  -2.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.96"
$ws.Range("E43").Value = "  +3.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.49"
$ws.Range("E44").Value = "  -4.23%  "
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "151.90"
$ws.Range("E46").Value = "  -2.87%  "
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.68"
$ws.Range("E49").Value = "  -2.33%  "
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.30"
$ws.Range("E51").Value = "  +2.03%  "
